$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $newValue)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.ClearFormats()
}

Set-TextValue "D2" "289.93"
Set-TextValue "E2" "-4.05%"
Set-TextValue "D3" "30.83"
Set-TextValue "E3" "-4.41%"
Set-TextValue "D4" "4.885"
Set-TextValue "E4" "-2.12%"
Set-TextValue "D5" "0.07190"
Set-TextValue "E5" "-9.04%"
Set-TextValue "D6" "7.675"
Set-TextValue "D7" "1.759"
Set-TextValue "E7" "-16.38%"
Set-TextValue "D8" "3.737"
Set-TextValue "E8" "-1.70%"
Set-TextValue "D9" "0.8946"
Set-TextValue "E9" "-3.58%"
Set-TextValue "D10" "0.1657"
Set-TextValue "E10" "-5.18%"
Set-TextValue "D11" "0.07490"
Set-TextValue "E11" "-5.12%"
Set-TextValue "D12" "0.08040"
Set-TextValue "E12" "-7.80%"
Set-TextValue "D13" "0.02987"
Set-TextValue "E13" "-4.88%"
Set-TextValue "D14" "0.09995"
Set-TextValue "E14" "-0.35%"
Set-TextValue "E15" "-1.31%"
Set-TextValue "D16" "0.005749"
Set-TextValue "E16" "-0.14%"
Set-TextValue "D18" "3.459"
Set-TextValue "E18" "-0.18%"
Set-TextValue "D19" "2.103"
Set-TextValue "E19" "-7.58%"
Set-TextValue "D20" "0.3278"
Set-TextValue "E20" "-0.28%"
Set-TextValue "D21" "0.1299"
Set-TextValue "E21" "0.64%"
Set-TextValue "D22" "4.410"
Set-TextValue "E22" "1.66%"
Set-TextValue "E23" "11.74%"
Set-TextValue "D24" "0.04467"
Set-TextValue "E24" "-3.07%"
Set-TextValue "D25" "0.001212"
Set-TextValue "E25" "-2.11%"
Set-TextValue "D26" "0.004021"
Set-TextValue "E26" "-9.91%"
Set-TextValue "E27" "0.07%"
Set-TextValue "D39" "0.01642"
Set-TextValue "E39" "-4.89%"
Set-TextValue "D40" "0.04332"
Set-TextValue "E40" "-9.66%"
Set-TextValue "D41" "0.007418"
Set-TextValue "E41" "-0.41%"
Set-TextValue "E42" "-3.73%"
Set-TextValue "D43" "0.002017"
Set-TextValue "E43" "-14.61%"
Set-TextValue "E44" "-0.97%"
Set-TextValue "E45" "-4.64%"
Set-TextValue "D46" "0.00000000751"
Set-TextValue "E46" "0.06%"
Set-TextValue "D47" "2.189"
Set-TextValue "E47" "166.77%"
Set-TextValue "D48" "0.003005"
Set-TextValue "E48" "-11.45%"
Set-TextValue "D49" "0.00002103"
Set-TextValue "E49" "0.06%"
Set-TextValue "D50" "0.0002003"
Set-TextValue "E50" "0.06%"
